# Update "想去人数" (wanted-to-go count) values on the "展览" and "全部类型" sheets
# to reflect newly generated output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1091
$ws1.Range("F4").Value = 1692
$ws1.Range("F5").Value = 761
$ws1.Range("F6").Value = 188

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1091
$ws4.Range("F4").Value = 1692
$ws4.Range("F6").Value = 761
$ws4.Range("F7").Value = 188
